# Incorporated WRI files sent 4/30/20
#
# 1) "About" sheet: insert a blank separator row and a new
#    "Comment (Brazil):" / note row near the top of the sheet.
# 2) "BPaFF-BITPTaP" sheet: flag hard coal and biomass as peakers
#    (their dependent formula cells recalc automatically).
# 3) "BPaFF-BDTPTPF" sheet: flag hard coal, hydro and biomass as
#    flexibility providers (their dependent formula cells recalc
#    automatically).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Push the old row 6 ("Notes") and everything below it down by two
# rows, opening up a blank row 5 and a new row 6.
$wsAbout.Rows.Item(5).Insert()
$wsAbout.Rows.Item(6).Insert()

$wsAbout.Cells.Item(6, 1).Value = "Comment (Brazil):"
$wsAbout.Cells.Item(6, 3).Value = "Hydro is not a peaker type, but can provide flexibility"

$wsAbout.Range("C6").Select()

# ---------------------------------------------------------------
# Sheet "BPaFF-BITPTaP"
# ---------------------------------------------------------------
$wsBITPTaP = $wb.Worksheets.Item("BPaFF-BITPTaP")

$wsBITPTaP.Cells.Item(2, 2).Value = 1   # hard coal
$wsBITPTaP.Cells.Item(9, 2).Value = 1   # biomass

# ---------------------------------------------------------------
# Sheet "BPaFF-BDTPTPF"
# ---------------------------------------------------------------
$wsBDTPTPF = $wb.Worksheets.Item("BPaFF-BDTPTPF")

$wsBDTPTPF.Cells.Item(2, 2).Value = 1   # hard coal
$wsBDTPTPF.Cells.Item(5, 2).Value = 1   # hydro
$wsBDTPTPF.Cells.Item(9, 2).Value = 1   # biomass
